$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Monitor/Alex Valenzuela-DellMonitor -> Computer/None-MacBookPro
$ws.Range("A2").Value = "Computer"
$ws.Range("B2").Value = "None-MacBookPro"
$ws.Range("D2").Value = "Apple Inc"
$ws.Range("E2").Value = "MacBook Pro"
$ws.Range("F2").Value = "SF2WN4N77MC"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "Laptop"

# Row 3: sebastian salgado-MacBookPro -> Sebastian Salgado-Latitude
$ws.Range("B3").Value = "Sebastian Salgado-Latitude"
$ws.Range("D3").Value = "Dell inc."
$ws.Range("E3").Value = "Latitude"
$ws.Range("F3").Value = "CS08BY3"
$ws.Range("M3").Value = "Sebastian Salgado"
